$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repull of data
$ws.Range("F2").Value = 5
$ws.Range("F4").Value = -1
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -5
$ws.Range("F10").Value = 4
$ws.Range("F12").Value = -6
$ws.Range("F13").Value = 10
$ws.Range("F14").Value = 2
$ws.Range("F15").Value = 0
